# LOM3013.docx content reshuffle
#
# The section headings (Objetivos, Docente(s) Responsável(eis), Programa
# resumido, Programa, Avaliação, Bibliografia, Requisitos) stay in place;
# only the paragraph(s) that hold the body text under each heading are
# swapped around. We target each body paragraph directly via the
# Paragraphs collection (1-based index) and scope Find/Replace to that
# paragraph's Range so that duplicate/overlapping text elsewhere in the
# document can never be matched by mistake.
#
# [char]11 is used for the manual line break (the same character Word's
# Find/Replace uses for "^l"), so that text that is supposed to stay
# within a single run/paragraph but spans multiple visual lines (joined
# by <w:br/>) can be matched / produced in one Find.Execute call.

$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceAll = 2
$brk = [char]11

# NOTE: this runtime's PowerShell does not bind named (-Param value)
# arguments to function param() blocks correctly, so every helper here
# is called with plain positional arguments instead.
function Replace-InParagraph($Index, $OldText, $NewText) {
    $p = $d.Paragraphs.Item($Index)
    $rng = $p.Range
    $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $NewText, $wdReplaceAll)
}

# --- "Objetivos" body (paragraph 6): was the long "Apresentação..." text,
#     becomes the short "Ligação atômica..." summary. ---
$p6Old = "Apresentação dos fundamentos da Ciência dos Materiais visando a introdução ao estudo das características microestruturais e das propriedades dos materiais, apresentação e discussão de exemplos práticos, bem como fornecer subsídios para o estudo das demais disciplinas do ciclo profissional."
$p6New = "Ligação atômica. Estrutura cristalina. Defeitos em cristais e estruturas não-cristalinas. Relação microestrutura-propriedade."
Replace-InParagraph 6 $p6Old $p6New

# --- "Docente(s) Responsável(eis)" body (paragraph 8): was the list of
#     four professors, becomes the objectives text + full program list +
#     the "método" paragraph text. ---
$p8Old = "6495737 - Durval Rodrigues Junior" + $brk + `
    "5983729 - Fernando Vernilli Junior" + $brk + `
    "984972 - Hugo Ricardo Zschommler Sandim" + $brk + `
    "7459752 - Maria Ismenia Sodero Toledo Faria"
$p8New = "Apresentação dos fundamentos da Ciência dos Materiais visando a introdução ao estudo das características microestruturais e das propriedades dos materiais, apresentação e discussão de exemplos práticos, bem como fornecer subsídios para o estudo das demais disciplinas do ciclo profissional." + $brk + `
    "1. Estrutura atômica, ligações atômicas: ligação iônica, ligação covalente, ligação metálica, ligação de Van der Waals, interações dipolo-dipolo e pontes de hidrogênio. Ligações atômicas e o coeficiente de expansão linear." + $brk + `
    "2. Estrutura cristalina: os sete sistemas e as quatorze redes de Bravais; estruturas de metais, cerâmicas e polímeros; direções e planos atômicos (notação de Miller), número de coordenação, empacotamento atômico linear a planar, Lei de Bragg e difração de raios-X;" + $brk + `
    "3. Defeitos em cristais e em estruturas amorfas: soluções sólidas (intersticiais e substitucionais); defeitos de ponto, defeitos de linha (discordâncias e sua dinâmica: movimentação e interação), defeitos bidimensionais (falhas de empilhamento, contornos de antifase, contornos de alto e de baixo ângulo), sólidos amorfos, vidros metálicos, defeitos tridimensionais (poros, trincas e inclusões)." + $brk + `
    "4. Relação microestrutura-propriedade: exemplos práticos e estudos de caso (propriedades mecânicas, elétricas e magnéticas)." + $brk + `
    "Em todos os itens, são abordados os aspectos práticos de cada tópico da ementa para ampliar as competências dos alunos, que serão trabalhados com Estudos de Caso. Viagem didática complementar" + $brk + `
    "Esta é uma disciplina de caráter fundamental, exigindo dedicação individual para assimilação das definições e conceitos. Isto envolve leitura concentrada para fixação dos conceitos teóricos e realização de exercícios numéricos. Duas provas escritas (P1 e P2) serão aplicadas e com pesos iguais. O desenvolvimento do aluno ao longo do curso será aferido e estimulado por meio de discussões sobre um dado tema, porém sem a atribuição de nota, por conta da subjetividade envolvida." + $brk + `
    ": A Nota final (NF) será calculada da seguinte maneira: NF = (0,4*P1 +0,4* P2+ 0,2*NT) / 3"
Replace-InParagraph 8 $p8Old $p8New

# --- "Programa resumido" body (paragraph 10): was the "Ligação
#     atômica..." summary, becomes the "recuperação" formula text. ---
$p10Old = "Ligação atômica. Estrutura cristalina. Defeitos em cristais e estruturas não-cristalinas. Relação microestrutura-propriedade."
$p10New = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: " + $brk + "MR = (NF + PR) / 2"
Replace-InParagraph 10 $p10Old $p10New

# --- "Programa" body (paragraph 12): was the four-item syllabus text,
#     becomes the nine-item bibliography list. ---
$p12Old = "1. Estrutura atômica, ligações atômicas: ligação iônica, ligação covalente, ligação metálica, ligação de Van der Waals, interações dipolo-dipolo e pontes de hidrogênio. Ligações atômicas e o coeficiente de expansão linear." + $brk + `
    "2. Estrutura cristalina: os sete sistemas e as quatorze redes de Bravais; estruturas de metais, cerâmicas e polímeros; direções e planos atômicos (notação de Miller), número de coordenação, empacotamento atômico linear a planar, Lei de Bragg e difração de raios-X;" + $brk + `
    "3. Defeitos em cristais e em estruturas amorfas: soluções sólidas (intersticiais e substitucionais); defeitos de ponto, defeitos de linha (discordâncias e sua dinâmica: movimentação e interação), defeitos bidimensionais (falhas de empilhamento, contornos de antifase, contornos de alto e de baixo ângulo), sólidos amorfos, vidros metálicos, defeitos tridimensionais (poros, trincas e inclusões)." + $brk + `
    "4. Relação microestrutura-propriedade: exemplos práticos e estudos de caso (propriedades mecânicas, elétricas e magnéticas)." + $brk + `
    "Em todos os itens, são abordados os aspectos práticos de cada tópico da ementa para ampliar as competências dos alunos, que serão trabalhados com Estudos de Caso. Viagem didática complementar"
$p12New = "1. CALLISTER Jr, W.D., RETHWISCH, D.G. Ciência e Engenharia de Materiais: Uma Introdução, 8ª ed., LTC Editora, 2013." + $brk + `
    "2. ASKELAND, D.R., PHULÉ, P.P., Ciência e Engenharia dos Materiais, CENGAGE, São Paulo, 2008." + $brk + `
    "3. SHACKELFORD, J.F., Ciência dos Materiais, 6a. ed., Pearson, 2008." + $brk + `
    "4. PADILHA, A.F., Materiais para Engenharia: Microestrutura e Propriedades, Hemus Editora, 1997." + $brk + `
    "5. PADILHA, A.F., Técnicas de Análise Microestrutural, Ed. Hemus, 1985." + $brk + `
    "6. REED-HILL, R.E., Princípios de Metalurgia Física, Guanabara Dois, 1982." + $brk + `
    "7. BRANDON, D.D., KAPLAN, W.D., Microstructural Characterization of Materials, 1st. ed., Wiley, 1999." + $brk + `
    "8. ASHBY, M.F., JONES, D.R.H., Engenharia de Materiais, Elsevier Editora, 2007." + $brk + `
    "9. ASHBY, M.F., SHERCLIFF, H., CEBON, D., Materials: Engineering, Science, Processing and Design, Butterworth-Heinemann, 2010."
Replace-InParagraph 12 $p12Old $p12New

# --- "Avaliação" body (paragraph 14): keeps its three bold labels
#     (Método:, Critério:, Norma de recuperação:) but the text after each
#     label is replaced with one of the professor names (each handled as
#     its own Find/Replace so the bold label runs are left untouched). ---
$p14aOld = "Esta é uma disciplina de caráter fundamental, exigindo dedicação individual para assimilação das definições e conceitos. Isto envolve leitura concentrada para fixação dos conceitos teóricos e realização de exercícios numéricos. Duas provas escritas (P1 e P2) serão aplicadas e com pesos iguais. O desenvolvimento do aluno ao longo do curso será aferido e estimulado por meio de discussões sobre um dado tema, porém sem a atribuição de nota, por conta da subjetividade envolvida."
Replace-InParagraph 14 $p14aOld "6495737 - Durval Rodrigues Junior"

$p14bOld = ": A Nota final (NF) será calculada da seguinte maneira: NF = (0,4*P1 +0,4* P2+ 0,2*NT) / 3"
Replace-InParagraph 14 $p14bOld "5983729 - Fernando Vernilli Junior"

$p14cOld = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: " + $brk + "MR = (NF + PR) / 2"
Replace-InParagraph 14 $p14cOld "984972 - Hugo Ricardo Zschommler Sandim"

# --- "Bibliografia" body (paragraph 16): was the nine-item bibliography
#     list, becomes the one remaining professor name. ---
$p16Old = "1. CALLISTER Jr, W.D., RETHWISCH, D.G. Ciência e Engenharia de Materiais: Uma Introdução, 8ª ed., LTC Editora, 2013." + $brk + `
    "2. ASKELAND, D.R., PHULÉ, P.P., Ciência e Engenharia dos Materiais, CENGAGE, São Paulo, 2008." + $brk + `
    "3. SHACKELFORD, J.F., Ciência dos Materiais, 6a. ed., Pearson, 2008." + $brk + `
    "4. PADILHA, A.F., Materiais para Engenharia: Microestrutura e Propriedades, Hemus Editora, 1997." + $brk + `
    "5. PADILHA, A.F., Técnicas de Análise Microestrutural, Ed. Hemus, 1985." + $brk + `
    "6. REED-HILL, R.E., Princípios de Metalurgia Física, Guanabara Dois, 1982." + $brk + `
    "7. BRANDON, D.D., KAPLAN, W.D., Microstructural Characterization of Materials, 1st. ed., Wiley, 1999." + $brk + `
    "8. ASHBY, M.F., JONES, D.R.H., Engenharia de Materiais, Elsevier Editora, 2007." + $brk + `
    "9. ASHBY, M.F., SHERCLIFF, H., CEBON, D., Materials: Engineering, Science, Processing and Design, Butterworth-Heinemann, 2010."
Replace-InParagraph 16 $p16Old "7459752 - Maria Ismenia Sodero Toledo Faria"

Write-Host "Content reshuffle complete."
